# feat: add 2022-Q1 data
#
# The existing "总计" (totals) sheet is renamed to "2022-Q1" and repurposed to
# hold the new quarter's per-fund detail rows (same shape as the other
# quarterly sheets). A brand-new "总计" sheet is inserted right after it,
# carrying the old totals-table content plus one new leading row summarising
# the 2022-Q1 quarter.

$wb = $excel.ActiveWorkbook

$total = $wb.Worksheets.Item("总计")

# Insert the new summary sheet immediately after the current "总计" sheet,
# before we rename anything (so sheet ordering matches: ... 2021-Q4, 2022-Q1, 总计).
$newTotal = $wb.Worksheets.Add($null, $total)
$newTotal.Name = "总计_tmp"

# Repurpose the old "总计" sheet as the new "2022-Q1" per-fund detail sheet.
$total.Name = "2022-Q1"
$newTotal.Name = "总计"

# ---------------------------------------------------------------------
# 2022-Q1 sheet: replace the old totals-table layout (A:D, 6 rows) with
# the per-fund detail layout (A:H, 3 rows) used by the other quarter tabs.
# ---------------------------------------------------------------------

# Drop the 3 surplus data rows left over from the old totals table.
$total.Rows.Item(6).EntireRow.Delete()
$total.Rows.Item(5).EntireRow.Delete()
$total.Rows.Item(4).EntireRow.Delete()

# Extend the header formatting (bold / bordered / centered style already on
# D1) across the new E:H header cells.
$total.Range("D1").Copy()
$total.Range("E1:H1").PasteSpecial(-4122)

$total.Cells.Item(1,2).Value = "基金代码"
$total.Cells.Item(1,3).Value = "基金名称"
$total.Cells.Item(1,4).Value = "基金规模"
$total.Cells.Item(1,5).Value = "股票总仓位"
$total.Cells.Item(1,6).Value = "仓位占比"
$total.Cells.Item(1,7).Value = "持有市值(亿元)"
$total.Cells.Item(1,8).Value = "仓位排名"

# Row 2: 006440
$total.Range("B2:D2").NumberFormat = "@"
$total.Cells.Item(2,2).Value = "006440"
$total.Cells.Item(2,3).Value = "中信建投中证500指数增强A"
$total.Cells.Item(2,4).Value = "5.78"
$total.Range("E2:G2").NumberFormat = "@"
$total.Cells.Item(2,5).Value = "94.71"
$total.Cells.Item(2,6).Value = "1.12"
$total.Cells.Item(2,7).Value = "0.0647"
$total.Cells.Item(2,8).Value = 3

# Row 3: 006441
$total.Range("B3:D3").NumberFormat = "@"
$total.Cells.Item(3,2).Value = "006441"
$total.Cells.Item(3,3).Value = "中信建投中证500指数增强C"
$total.Cells.Item(3,4).Value = "3.11"
$total.Range("E3:G3").NumberFormat = "@"
$total.Cells.Item(3,5).Value = "94.71"
$total.Cells.Item(3,6).Value = "1.12"
$total.Cells.Item(3,7).Value = "0.0348"
$total.Cells.Item(3,8).Value = 3

# ---------------------------------------------------------------------
# 总计 sheet: rebuild the totals table (A:D, 7 rows) - the original 6 rows
# of history plus one new leading row for 2022-Q1 - using the same cell
# styling as the other quarter tabs' header / index-column cells.
# ---------------------------------------------------------------------

$refSheet = $wb.Worksheets.Item("2021-Q4")

# Body-row style (bold/bordered index column A, plain B:D) for all 6 data rows.
$refSheet.Range("A2:D2").Copy()
$newTotal.Range("A2:D7").PasteSpecial(-4122)

# Header-row style.
$refSheet.Range("B1:D1").Copy()
$newTotal.Range("B1:D1").PasteSpecial(-4122)

$newTotal.Cells.Item(1,2).Value = "日期"
$newTotal.Cells.Item(1,3).Value = "持有数量(只)"
$newTotal.Cells.Item(1,4).Value = "持有市值(亿元)"

$totalsData = @(
  @(0, "2022-Q1", 2,  0.1),
  @(1, "2021-Q4", 7,  0.6),
  @(2, "2021-Q3", 7,  1.11),
  @(3, "2021-Q2", 13, 4.73),
  @(4, "2021-Q1", 15, 3.68),
  @(5, "2020-Q4", 3,  0.47)
)

$r = 2
foreach ($row in $totalsData) {
    $newTotal.Cells.Item($r,1).Value = $row[0]
    $newTotal.Cells.Item($r,2).Value = $row[1]
    $newTotal.Cells.Item($r,3).Value = $row[2]
    $newTotal.Cells.Item($r,4).Value = $row[3]
    $r++
}
